$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "level" and "level_lab" columns (G, H); column F gets
# overwritten in place below (its header/content changes from
# "definition" to "level").
$ws.Range("G:H").Delete() | Out-Null

# New header row
$ws.Range("A1").Value = "varName"
$ws.Range("B1").Value = "label"
$ws.Range("C1").Value = "shortName"
$ws.Range("D1").Value = "Levels"
$ws.Range("E1").Value = "Labels"
$ws.Range("F1").Value = "level"

# Data rows (data in new order, including 5 brand-new rows)
$data = @(
    @("fhh", "Gender of Household Head", "Male/Female Household Head", "0,1", "Male,Female", "All"),
    @("covid_shock", "Household Impacted by COVID", "COVID Impact", "1,2", "Yes,No", "All"),
    @("ag_comm", "Membership in an Agricultural Community", "Association Membership", "1,2", "Yes,No", "All"),
    @("livestock_area", "Land Area for Livestock", "Farm Size", "0,1,2,3", "0 ha,>0-2 ha,>2-4 ha,>4 ha", "All"),
    @("feed_Chickens", "Chicken Feeding Practice", "Type of Feed", "1,2,3", "Only Grazing,Mixed,Only Feed", "Poultry"),
    @("drought", "Household Experienced Drought", "Drought", "0,1", "No Drought,Drought", "All"),
    @("flood", "Household Experienced Flood", "Flooding", "0,1", "No Flooding,Flooding", "All"),
    @("irrigation", "Household Used Irrigation", "Irrigation Use", "0,1", "No Irrigation,Irrigation", "All"),
    @("feed_Cattle", "Livestock Feeding Practice", "Type of Feed", "1,2,3", "Only Grazing,Mixed,Only Feed", "Livestock")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# Column widths - recompute "best fit" now that columns A and B contain
# longer strings than before.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null

# Selection moves to A11 to match post-edit state
$ws.Range("A11").Select() | Out-Null
